$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.819.02"
$ws.Range("E2").Value = "  -0.67%  "

# Row 3
$ws.Range("D3").Value = "2.346.07"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.73"
$ws.Range("E5").Value = "  -0.85%  "

# Row 6
$ws.Range("E6").Value = "  -4.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.38"
$ws.Range("E7").Value = "  -4.65%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  -4.54%  "

# Row 10
$ws.Range("E10").Value = "  -0.92%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.54"
$ws.Range("E11").Value = "  +3.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.78"
$ws.Range("E12").Value = "  -2.14%  "

# Row 13
$ws.Range("E13").Value = "  -0.60%  "

# Row 14
$ws.Range("E14").Value = "  -3.36%  "

# Row 15
$ws.Range("D15").Value = "2.696.23"
$ws.Range("E15").Value = "  -0.94%  "

# Row 16
$ws.Range("E16").Value = "  -3.82%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.903"
$ws.Range("E17").Value = "  -3.13%  "

# Row 18
$ws.Range("D18").Value = "2.343.84"
$ws.Range("E18").Value = "  -0.91%  "

# Row 19
$ws.Range("D19").Value = "43.755.32"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  +0.20%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "78.37"
$ws.Range("E21").Value = "  +0.63%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.66"
$ws.Range("E22").Value = "  -0.83%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.13"
$ws.Range("E23").Value = "  -3.02%  "

# Row 24
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.81"
$ws.Range("E24").Value = "  +1.79%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.12%  "

# Row 26
$ws.Range("E26").Value = "  +2.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("E27").Value = "  -1.97%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.48"
$ws.Range("E28").Value = "  +10.88%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.42"
$ws.Range("E29").Value = "  -4.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.64"
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("E31").Value = "  -4.50%  "

# Row 32
$ws.Range("E32").Value = "  -0.62%  "

# Row 33
$ws.Range("E33").Value = "  -2.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0745"
$ws.Range("E34").Value = "  -2.87%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.07"
$ws.Range("E35").Value = "  -6.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.35"
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
$ws.Range("E38").Value = "  -0.98%  "

# Row 39
$ws.Range("E39").Value = "  -2.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.60"
$ws.Range("E40").Value = "  +14.56%  "

# Row 41
$ws.Range("E41").Value = "  -4.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.83"
$ws.Range("E42").Value = "  +15.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.17"
$ws.Range("E43").Value = "  -0.38%  "

# Row 44
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.105"
$ws.Range("E44").Value = "  -2.16%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.77"
$ws.Range("E45").Value = "  -3.21%  "

# Row 46
$ws.Range("E46").Value = "  -11.55%  "

# Row 47
$ws.Range("E47").Value = "  +0.11%  "

# Row 48
$ws.Range("E48").Value = "  -3.24%  "

# Row 49
$ws.Range("E49").Value = "  -4.07%  "

# Row 50
$ws.Range("E50").Value = "  -5.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.18"
$ws.Range("E51").Value = "  -4.41%  "
